# Memory map updates: expand SRAM from 64KB to 128KB window,
# select E10 cell on "Memory Map" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Memory Map")

# A7: (64-10)*1024 -> (128)*1024
$ws.Range("A7").Formula = "=(128)*1024"

# A8: literal 0 -> (128-10)*1024 formula
$ws.Range("A8").Formula = "=(128-10)*1024"

# Update the selection to E10 on the Memory Map sheet
$ws.Activate()
$ws.Range("E10").Select()
